$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new value in D16 ("yellow" wire color, matching C16's sabertooth rx row)
$ws.Range("D16").Value = "yellow"

# Update the active cell / selection to D16, matching the recorded cursor move
$ws.Range("D16").Select()
